$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 605; this shifts existing rows 605-685 down to 606-686
# (Excel copies the cell formatting from the row above automatically).
$ws.Rows.Item(605).Insert()

# Fill the newly inserted row 605 with the new weekly record.
$ws.Range("A605").Value2 = 8
$ws.Range("B605").Value2 = "Terminal La Palmera de La Serena"
$ws.Range("C605").Value2 = "Coquimbo"
$ws.Range("D605").Value2 = 45142
$ws.Range("D605").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E605").Value2 = 4
$ws.Range("F605").Value2 = 100112017
$ws.Range("G605").Value2 = "Apio"
$ws.Range("H605").Value2 = "Americana (o)"
$ws.Range("I605").Value2 = "Primera"
$ws.Range("J605").Value2 = 1200
$ws.Range("K605").Value2 = 6000
$ws.Range("L605").Value2 = 7000
$ws.Range("M605").Value2 = 6500
$ws.Range("N605").Value2 = "`$/docena de matas"
$ws.Range("O605").Value2 = "Provincia del Elquí"
$ws.Range("P605").Value2 = 1083
$ws.Range("Q605").Value2 = 6
$ws.Range("R605").Value2 = "Hortaliza"
